# Atualização de bases das ligas, do dia: 17-02-2024 às 11:11
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 69 (swap with row 70's original values) ---
$ws.Cells.Item(69, 2).Value  = 5240690
$ws.Cells.Item(69, 6).Value  = "FCSB"
$ws.Cells.Item(69, 7).Value  = "ACS Sepsi"
$ws.Cells.Item(69, 11).Value = 2.05
$ws.Cells.Item(69, 12).Value = 3.25
$ws.Cells.Item(69, 13).Value = 3.4
$ws.Cells.Item(69, 14).Value = 1.833
$ws.Cells.Item(69, 15).Value = 3.5
$ws.Cells.Item(69, 16).Value = 3.8
$ws.Cells.Item(69, 17).Value = -0.5
$ws.Cells.Item(69, 18).Value = 1.85
$ws.Cells.Item(69, 19).Value = 2
$ws.Cells.Item(69, 20).Value = 2.5
$ws.Cells.Item(69, 21).Value = 1.95
$ws.Cells.Item(69, 22).Value = 1.9
$ws.Cells.Item(69, 23).Value = 0.833
$ws.Cells.Item(69, 26).Value = 0.8500000000000001
$ws.Cells.Item(69, 29).Value = 0.8999999999999999

# --- Row 70 (swap with row 69's original values) ---
$ws.Cells.Item(70, 2).Value  = 5240692
$ws.Cells.Item(70, 6).Value  = "AFC Hermannstadt"
$ws.Cells.Item(70, 7).Value  = "FC U Craiova 1948"
$ws.Cells.Item(70, 11).Value = 3.2
$ws.Cells.Item(70, 12).Value = 3
$ws.Cells.Item(70, 13).Value = 2.2
$ws.Cells.Item(70, 14).Value = 2.875
$ws.Cells.Item(70, 15).Value = 2.875
$ws.Cells.Item(70, 16).Value = 2.45
$ws.Cells.Item(70, 17).Value = 0
$ws.Cells.Item(70, 18).Value = 2.1
$ws.Cells.Item(70, 19).Value = 1.775
$ws.Cells.Item(70, 20).Value = 2
$ws.Cells.Item(70, 21).Value = 2
$ws.Cells.Item(70, 22).Value = 1.85
$ws.Cells.Item(70, 23).Value = 1.875
$ws.Cells.Item(70, 26).Value = 1.1
$ws.Cells.Item(70, 29).Value = 0.8500000000000001

# --- Row 353 ---
$ws.Cells.Item(353, 2).Value  = 6836256
$ws.Cells.Item(353, 5).Value  = 45339.6875
$ws.Cells.Item(353, 6).Value  = "FC U Craiova 1948"
$ws.Cells.Item(353, 7).Value  = "CFR Cluj"
$ws.Cells.Item(353, 11).Value = 3.1
$ws.Cells.Item(353, 12).Value = 3.1
$ws.Cells.Item(353, 13).Value = 2.25
$ws.Cells.Item(353, 14).Value = 3.5
$ws.Cells.Item(353, 15).Value = 3.4
$ws.Cells.Item(353, 16).Value = 1.95
$ws.Cells.Item(353, 17).Value = 0.5
$ws.Cells.Item(353, 18).Value = 1.875
$ws.Cells.Item(353, 19).Value = 1.975
$ws.Cells.Item(353, 20).Value = 2.5
$ws.Cells.Item(353, 21).Value = 1.9
$ws.Cells.Item(353, 22).Value = 1.95

# --- Row 354 ---
$ws.Cells.Item(354, 2).Value  = 6836255
$ws.Cells.Item(354, 5).Value  = 45340.51041666666
$ws.Cells.Item(354, 6).Value  = "FC Botosani"
$ws.Cells.Item(354, 7).Value  = "CS U Craiova"
$ws.Cells.Item(354, 11).Value = 3.4
$ws.Cells.Item(354, 12).Value = 3.4
$ws.Cells.Item(354, 13).Value = 2
$ws.Cells.Item(354, 14).Value = 4
$ws.Cells.Item(354, 15).Value = 3.5
$ws.Cells.Item(354, 16).Value = 1.8
$ws.Cells.Item(354, 17).Value = 0.5
$ws.Cells.Item(354, 18).Value = 2
$ws.Cells.Item(354, 19).Value = 1.85
$ws.Cells.Item(354, 20).Value = 2.25
$ws.Cells.Item(354, 21).Value = 1.825
$ws.Cells.Item(354, 22).Value = 2.025

# --- Row 355 ---
$ws.Cells.Item(355, 2).Value  = 6836258
$ws.Cells.Item(355, 5).Value  = 45340.625
$ws.Cells.Item(355, 6).Value  = "Dinamo Bucharest"
$ws.Cells.Item(355, 7).Value  = "Otelul Galati"
$ws.Cells.Item(355, 11).Value = 2.5
$ws.Cells.Item(355, 12).Value = 2.9
$ws.Cells.Item(355, 13).Value = 2.9
$ws.Cells.Item(355, 14).Value = 2.5
$ws.Cells.Item(355, 15).Value = 2.9
$ws.Cells.Item(355, 16).Value = 2.9
$ws.Cells.Item(355, 17).Value = 0
$ws.Cells.Item(355, 18).Value = 1.75
$ws.Cells.Item(355, 19).Value = 2.125
$ws.Cells.Item(355, 20).Value = 2
$ws.Cells.Item(355, 21).Value = 2.05
$ws.Cells.Item(355, 22).Value = 1.8

# --- Row 356 ---
$ws.Cells.Item(356, 2).Value  = 6836259
$ws.Cells.Item(356, 5).Value  = 45341.52083333334
$ws.Cells.Item(356, 6).Value  = "ACS Sepsi"
$ws.Cells.Item(356, 7).Value  = "Farul Constanta"
$ws.Cells.Item(356, 11).Value = 2.3
$ws.Cells.Item(356, 12).Value = 3.2
$ws.Cells.Item(356, 13).Value = 2.9
$ws.Cells.Item(356, 14).Value = 2.25
$ws.Cells.Item(356, 15).Value = 3.2
$ws.Cells.Item(356, 16).Value = 3
$ws.Cells.Item(356, 17).Value = -0.25
$ws.Cells.Item(356, 18).Value = 2
$ws.Cells.Item(356, 19).Value = 1.85
$ws.Cells.Item(356, 20).Value = 2.5
$ws.Cells.Item(356, 21).Value = 2.05
$ws.Cells.Item(356, 22).Value = 1.8

# --- Row 357 ---
$ws.Cells.Item(357, 2).Value  = 6836257
$ws.Cells.Item(357, 5).Value  = 45341.64583333334
$ws.Cells.Item(357, 6).Value  = "Universitatea Cluj"
$ws.Cells.Item(357, 7).Value  = "FCSB"
$ws.Cells.Item(357, 11).Value = 3.6
$ws.Cells.Item(357, 12).Value = 3.2
$ws.Cells.Item(357, 13).Value = 2
$ws.Cells.Item(357, 14).Value = 3.6
$ws.Cells.Item(357, 15).Value = 3.2
$ws.Cells.Item(357, 16).Value = 2.05
$ws.Cells.Item(357, 17).Value = 0.5
$ws.Cells.Item(357, 18).Value = 1.8
$ws.Cells.Item(357, 19).Value = 2.05
$ws.Cells.Item(357, 20).Value = 2.5
$ws.Cells.Item(357, 21).Value = 2.025
$ws.Cells.Item(357, 22).Value = 1.825

# --- Row 358 no longer exists: delete it, which also fixes the sheet dimension ---
$ws.Rows(358).Delete()
